function Replace-ParaText($paraRange, $oldText, $newText) {
    $found = $paraRange.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, "@", 2)
    if (-not $found) { Write-Host "NOT FOUND: $oldText" }
    # after the replace, $paraRange collapses onto the inserted '@' seed character;
    # expand over it, insert the real text right after it (inheriting its run
    # formatting), then delete the one-character seed so only the new text remains.
    $paraRange.Collapse(1)
    $paraRange.MoveEnd(1, 1)
    $paraRange.InsertAfter($newText)
    $seed = $paraRange.Document.Range($paraRange.Start, $paraRange.Start + 1)
    $seed.Delete()
}

$d = $word.ActiveDocument

# --- Item 1: "Efetuar login na aplicação;" -> "Verificar o calendário final de Ti do ano letivo 2019/2020;"
Replace-ParaText $d.Paragraphs.Item(24).Range "Efetuar login na aplicação;" "Verificar o calendário final de Ti do ano letivo 2019/2020;"

# --- Item 2: "Proceder à criação de um novo calendário;" -> "Importar ficheiro .csv;"
Replace-ParaText $d.Paragraphs.Item(25).Range "Proceder à criação de um novo calendário;" "Importar ficheiro .csv;"

# --- Item 3: "Marcar um exame no calendário acabado de criar;" -> "Verificar individualmente quantas disciplinas, salas e docentes existem;"
Replace-ParaText $d.Paragraphs.Item(26).Range "Marcar um exame no calendário acabado de criar;" "Verificar individualmente quantas disciplinas, salas e docentes existem;"

# --- Item 4: "Abrir outro calendário através do menu;" -> 'Pesquisar por "Ti" na barra de pesquisa e abrir o calendário "Ti - 1º Ano - 1º Semestre";'
Replace-ParaText $d.Paragraphs.Item(27).Range "Abrir outro calendário através do menu;" 'Pesquisar por "Ti" na barra de pesquisa e abrir o calendário "Ti - 1º Ano - 1º Semestre";'

# --- Item 5: 'Realizar a pesquisa de um calendário recorrendo à caixa de pesquisa existente no “Calendários”;' -> "Criar um novo calendário para o curso de Ti;"
Replace-ParaText $d.Paragraphs.Item(28).Range "Realizar a pesquisa de um calendário recorrendo à caixa de pesquisa existente no “Calendários”;" "Criar um novo calendário para o curso de Ti;"

# --- Item 6: "Utilizar os filtros de pesquisa de curso, época e semestre;" -> "Mova matemática para o período da manhã do dia 14;"
Replace-ParaText $d.Paragraphs.Item(29).Range "Utilizar os filtros de pesquisa de curso, época e semestre;" "Mova matemática para o período da manhã do dia 14;"

# --- Item 7: "Importar um .csv de disciplinas;" -> 'Colocar "Segurança Inf." num período da noite;'
Replace-ParaText $d.Paragraphs.Item(30).Range "Importar um .csv de disciplinas;" 'Colocar "Segurança Inf." num período da noite;'

# --- Item 8: "Fazer a configuração manual dos dados de uma sala de aulas;" -> "Exportar para um .pdf;"
Replace-ParaText $d.Paragraphs.Item(31).Range "Fazer a configuração manual dos dados de uma sala de aulas;" "Exportar para um .pdf;"

# --- Item 9: "Exportar os dados para .pdf;" -> "Fazer log out."
Replace-ParaText $d.Paragraphs.Item(32).Range "Exportar os dados para .pdf;" "Fazer log out."

# --- Item 10: paragraph "Realizar Log Out da aplicação." is removed entirely (merged into item 9)
$d.Paragraphs.Item(33).Range.Delete() | Out-Null

# --- Bookmark "_Hlk88835671" spans from the very start of item 1's paragraph to
#     the end of item 9's paragraph text (right after "Fazer log out.").
$startPos = $d.Paragraphs.Item(24).Range.Start
$endRange = $d.Paragraphs.Item(32).Range
$endRange.Collapse(0)
$endRange.MoveEnd(1, -1)   # move end back before the paragraph mark
$bmRange = $d.Range($startPos, $endRange.End)
$d.Bookmarks.Add("_Hlk88835671", $bmRange) | Out-Null

Write-Host "done"
